$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A134").Value = 133
$ws.Range("B134").Value = 1
$ws.Range("C134").Value = "2024-06-17 18:16:51"
$ws.Range("D134").Value = 200
$ws.Range("E134").Value = 21

$ws.Range("A135").Value = 134
$ws.Range("B135").Value = 2
$ws.Range("C135").Value = "2024-06-17 18:16:51"
$ws.Range("D135").Value = 200
$ws.Range("E135").Value = 2
